$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.645.22"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.502.02"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'573.62"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'166.24"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.512"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").Value = "2.500.05"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'0.168"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "'4.93"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "2.959.70"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "69.591.78"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "'0.0000176"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "'24.68"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "2.497.56"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'11.18"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "'7.50"
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").Value = "'348.61"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'1.93"
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'70.77"
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "'8.71"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "2.629.63"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").Value = "0.0₃0887"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'456.85"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "  -5.94%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'157.29"
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("D37").Value = "'0.115"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "'18.33"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "'0.316"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("D42").Value = "'4.68"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'2.19"
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("E46").Value = "  -8.18%  "
$ws.Range("D47").Value = "'140.96"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "'3.47"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").Value = "'0.517"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").Value = "  -0.49%  "
